# Complete rebuild of database (db name: cmms2)
# Insert two new rows above the existing "COMMIT TRANSACTION;" row and
# populate them with the new ALTER SEQUENCE statements, wrapping the text
# and widening column A to fit the new (longer) content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing row down by inserting two new rows at the top.
$ws.Rows("1:2").Insert()

# Populate the two new rows with the new SQL statements.
$ws.Range("A1").Value = "ALTER SEQUENCE orders_order_id_seq RESTART WITH 10;"
$ws.Range("A2").Value = "ALTER SEQUENCE persons_person_id_seq RESTART WITH 10;"

# Wrap the text in the two new cells (creates/applies the new cell style).
$ws.Range("A1:A2").WrapText = $true

# Widen column A so the longer statements are readable.
$ws.Columns("A:A").ColumnWidth = 64.83
